$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours for "Purpose, scope, definitions" (row 3)
$ws.Range("B3").Value = 2.5

# Update hours for "Functional requirements" (row 8)
$ws.Range("B8").Value = 9.5

# Add new hours entry for "Formal analysis using Alloy" (row 10)
$ws.Range("B10").Value = 0.5

# Update the selected cell in the sheet view
$ws.Range("B3").Select()
